$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.96799999999998
$ws.Range("A6").Value = -22.79870000000001
$ws.Range("A7").Value = -21.85660000000001
$ws.Range("B7").Value = 5.355100000000004
$ws.Range("B12").Value = 4.421200000000001
$ws.Range("E12").Value = 17.88640000000002
$ws.Range("C13").Value = -13.78529999999998
$ws.Range("C14").Value = -14.18009999999999
$ws.Range("B15").Value = 4.848399999999995
$ws.Range("A16").Value = -21.56209999999998
$ws.Range("C16").Value = -12.1471
$ws.Range("C19").Value = -12.297
$ws.Range("A20").Value = -22.83640000000001
$ws.Range("B20").Value = 4.698599999999995
$ws.Range("B21").Value = 10.3548
$ws.Range("B22").Value = 10.47060000000001
$ws.Range("C22").Value = -12.46480000000001
$ws.Range("E22").Value = 16.94400000000001
$ws.Range("B23").Value = 9.601900000000004
$ws.Range("A28").Value = -22.13539999999999
$ws.Range("A29").Value = -21.67979999999999
$ws.Range("B29").Value = 5.478000000000002
$ws.Range("E29").Value = 17.23990000000001
$ws.Range("A32").Value = -21.33620000000001
$ws.Range("B34").Value = 9.673000000000007
$ws.Range("E34").Value = 17.24020000000001
$ws.Range("C36").Value = -12.73420000000001
$ws.Range("A40").Value = -19.5839
$ws.Range("B42").Value = 9.512599999999994
$ws.Range("B43").Value = 6.176400000000004
$ws.Range("E43").Value = 17.37200000000001
$ws.Range("B44").Value = 4.610900000000004
$ws.Range("B45").Value = 5.210900000000001
$ws.Range("A46").Value = -22.3149
$ws.Range("B46").Value = 5.481399999999997
$ws.Range("C46").Value = -13.11599999999999
$ws.Range("E48").Value = 17.4487
$ws.Range("B50").Value = 4.772699999999997
$ws.Range("C50").Value = -13.73749999999999
$ws.Range("A51").Value = -22.22599999999999
$ws.Range("B51").Value = 5.535099999999997
$ws.Range("A52").Value = -22.12209999999999
$ws.Range("A57").Value = -22.82690000000001
$ws.Range("A59").Value = -22.0704
$ws.Range("E60").Value = 15.82700000000001
$ws.Range("A62").Value = -22.10940000000001
$ws.Range("A66").Value = -21.46660000000001
$ws.Range("B66").Value = 4.714199999999997
$ws.Range("B67").Value = 4.916399999999999
$ws.Range("E68").Value = 17.89980000000001
$ws.Range("E70").Value = 17.97770000000002
$ws.Range("A73").Value = -20.1711
$ws.Range("E73").Value = 17.26510000000001
$ws.Range("A74").Value = -22.02219999999999
$ws.Range("B79").Value = 9.717800000000006
$ws.Range("B84").Value = 5.624999999999997
$ws.Range("E87").Value = 16.32859999999999
$ws.Range("A92").Value = -21.40010000000002
$ws.Range("B92").Value = 4.857299999999998
$ws.Range("E92").Value = 18.77770000000002
$ws.Range("C95").Value = -11.56980000000001
$ws.Range("B97").Value = 5.713700000000001
$ws.Range("C97").Value = -11.4454
$ws.Range("A100").Value = -22.2135
$ws.Range("E101").Value = 16.83900000000001
